$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts all existing data rows down by one),
# picking up the same (unstyled) formatting as the data row that follows.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:D2").ClearFormats()

# Populate the new row 2 with the latest date and the same price values as the rest.
# Force the date cell to be stored as plain text (matching the other date cells),
# not auto-converted to a date serial number.
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2025-12-09"
$ws.Cells.Item(2, 1).Style = "Normal"

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
